$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '39.556.04'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.154.39'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.92'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.625'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.97'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.390'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0842'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.88'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.473.16'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.74'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.803'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.47'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.165.37'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '39.560.32'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.47'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.06'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.98'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.71'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.78'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.68'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'THORChain'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.89'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.30%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0614'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.67'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.40'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.11'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +22.81%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.55'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.66'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.515.31'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.82'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0915'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '49.63'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +8.07%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.11%  '
